$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 8666.923000000001
$ws.Range("I62").Value = 9312.424000000001
$ws.Range("K62").Value = 9312.424000000001
$ws.Range("M62").Value = -8688.424000000001
$ws.Range("H65").Value = 8666.923000000001
$ws.Range("I65").Value = 9312.424000000001
$ws.Range("K65").Value = 46562.12
$ws.Range("M65").Value = -43442.12
$ws.Range("H106").Value = 19679.79
$ws.Range("I106").Value = 22963.125
$ws.Range("J106").Value = 2168.6667
$ws.Range("K106").Value = 22963.125
$ws.Range("L106").Value = 2168.6667
$ws.Range("M106").Value = -22332.125
$ws.Range("N106").Value = -3430.6667
$ws.Range("H107").Value = 467.82608
$ws.Range("I107").Value = 510.6316
$ws.Range("J107").Value = 264.5
$ws.Range("K107").Value = 510.6316
$ws.Range("L107").Value = 264.5
$ws.Range("M107").Value = 1409.3684
$ws.Range("N107").Value = -4104.5
$ws.Range("H121").Value = 1771.8572
$ws.Range("J121").Value = 1771.8572
$ws.Range("L121").Value = 5315.571599999999
$ws.Range("N121").Value = -8809.571599999999
$ws.Range("H131").Value = 4502.6113
$ws.Range("I131").Value = 616.4167
$ws.Range("J131").Value = 12275
$ws.Range("K131").Value = 1849.2501
$ws.Range("L131").Value = 36825
$ws.Range("M131").Value = 3190.7499
$ws.Range("N131").Value = -46905
$ws.Range("H137").Value = 51332.6
$ws.Range("I137").Value = 112381.555
$ws.Range("J137").Value = 1383.4546
$ws.Range("K137").Value = 337144.665
$ws.Range("L137").Value = 4150.3638
$ws.Range("M137").Value = -334594.665
$ws.Range("N137").Value = -9250.363799999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7203.6665
$ws.Range("I2").Value = 654.5454999999999
$ws.Range("K2").Value = 654.5454999999999
$ws.Range("M2").Value = -541.5454999999999
$ws.Range("H45").Value = 1430.3
$ws.Range("I45").Value = 1499.5
$ws.Range("J45").Value = 1153.5
$ws.Range("K45").Value = 1499.5
$ws.Range("L45").Value = 1153.5
$ws.Range("M45").Value = -1122.5
$ws.Range("N45").Value = -1907.5
$ws.Range("H61").Value = 1490.2609
$ws.Range("I61").Value = 1253.8
$ws.Range("K61").Value = 1253.8
$ws.Range("M61").Value = -1041.8
$ws.Range("H74").Value = 80360.234
$ws.Range("I74").Value = 103762.4
$ws.Range("J74").Value = 2353
$ws.Range("K74").Value = 103762.4
$ws.Range("L74").Value = 2353
$ws.Range("M74").Value = -102888.4
$ws.Range("N74").Value = -4101
$ws.Range("H77").Value = 80360.234
$ws.Range("I77").Value = 103762.4
$ws.Range("J77").Value = 2353
$ws.Range("K77").Value = 518812
$ws.Range("L77").Value = 11765
$ws.Range("M77").Value = -514444
$ws.Range("N77").Value = -20501
$ws.Range("H116").Value = 7203.6665
$ws.Range("I116").Value = 654.5454999999999
$ws.Range("K116").Value = 654.5454999999999
$ws.Range("M116").Value = 1639.4545
$ws.Range("H132").Value = 1968130.2
$ws.Range("I132").Value = 2270075.8
$ws.Range("J132").Value = 835834.7
$ws.Range("K132").Value = 6810227.399999999
$ws.Range("L132").Value = 2507504.1
$ws.Range("M132").Value = -6807697.399999999
$ws.Range("N132").Value = -2512564.1
$ws.Range("H136").Value = 1490.2609
$ws.Range("I136").Value = 1253.8
$ws.Range("K136").Value = 3761.4
$ws.Range("M136").Value = -1211.4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7203.6665
$ws.Range("I3").Value = 654.5454999999999
$ws.Range("K3").Value = 654.5454999999999
$ws.Range("M3").Value = -540.5454999999999
$ws.Range("H128").Value = 900
$ws.Range("I128").Value = 900
$ws.Range("K128").Value = 2700
$ws.Range("M128").Value = -210
$ws.Range("H134").Value = 44638.81
$ws.Range("I134").Value = 2083.2104
$ws.Range("J134").Value = 160146.86
$ws.Range("K134").Value = 6249.6312
$ws.Range("L134").Value = 480440.58
$ws.Range("M134").Value = -3714.6312
$ws.Range("N134").Value = -485510.58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 875
$ws.Range("I16").Value = 810
$ws.Range("K16").Value = 810
$ws.Range("M16").Value = -523
$ws.Range("H58").Value = 2026.2
$ws.Range("I58").Value = 1637.3334
$ws.Range("J58").Value = 2192.8572
$ws.Range("K58").Value = 1637.3334
$ws.Range("L58").Value = 2192.8572
$ws.Range("M58").Value = -1434.3334
$ws.Range("N58").Value = -2598.8572
$ws.Range("H107").Value = 460.78125
$ws.Range("I107").Value = 349.88235
$ws.Range("K107").Value = 349.88235
$ws.Range("M107").Value = 1570.11765
$ws.Range("H113").Value = 875
$ws.Range("I113").Value = 810
$ws.Range("K113").Value = 810
$ws.Range("M113").Value = 1360
$ws.Range("H122").Value = 862.129
$ws.Range("I122").Value = 736.875
$ws.Range("J122").Value = 1291.5714
$ws.Range("K122").Value = 2210.625
$ws.Range("L122").Value = 3874.7142
$ws.Range("M122").Value = 239.375
$ws.Range("N122").Value = -8774.7142
$ws.Range("H132").Value = 102934.4
$ws.Range("I132").Value = 334908
$ws.Range("J132").Value = 3517.1428
$ws.Range("K132").Value = 1004724
$ws.Range("L132").Value = 10551.4284
$ws.Range("M132").Value = -1002194
$ws.Range("N132").Value = -15611.4284
$ws.Range("H136").Value = 2026.2
$ws.Range("I136").Value = 1637.3334
$ws.Range("J136").Value = 2192.8572
$ws.Range("K136").Value = 4912.0002
$ws.Range("L136").Value = 6578.571599999999
$ws.Range("M136").Value = -2362.0002
$ws.Range("N136").Value = -11678.5716
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 684.05884
$ws.Range("I122").Value = 520
$ws.Range("J122").Value = 734.53845
$ws.Range("K122").Value = 4680
$ws.Range("L122").Value = 6610.84605
$ws.Range("M122").Value = -2230
$ws.Range("N122").Value = -11510.84605
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 74160.07000000001
$ws.Range("I132").Value = 1450
$ws.Range("K132").Value = 4350
$ws.Range("M132").Value = -1820
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 434300.25
$ws.Range("I132").Value = 205671.5
$ws.Range("K132").Value = 617014.5
$ws.Range("M132").Value = -614484.5
$ws.Range("H136").Value = 386903.38
$ws.Range("I136").Value = 771126.75
$ws.Range("J136").Value = 2680
$ws.Range("K136").Value = 2313380.25
$ws.Range("L136").Value = 8040
$ws.Range("M136").Value = -2310830.25
$ws.Range("N136").Value = -13140
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 201
$ws.Range("I107").Value = 176
$ws.Range("K107").Value = 528
$ws.Range("M107").Value = 1392
$ws.Range("H113").Value = 212.0625
$ws.Range("I113").Value = 195.86957
$ws.Range("K113").Value = 587.60871
$ws.Range("M113").Value = 1582.39129
$ws.Range("H132").Value = 12484.1
$ws.Range("I132").Value = 1580
$ws.Range("J132").Value = 17157.285
$ws.Range("K132").Value = 4740
$ws.Range("L132").Value = 51471.855
$ws.Range("M132").Value = -2210
$ws.Range("N132").Value = -56531.855
$ws.Range("H136").Value = 4287456
$ws.Range("I136").Value = 7144425
$ws.Range("K136").Value = 21433275
$ws.Range("M136").Value = -21430725
